$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B-E) - subject/set identifiers updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (columns B-E) - CON meanEMG legmaxROM updated
$ws.Range("B2").Value = 87.672128106644124
$ws.Range("C2").Value = 60.750259275712338
$ws.Range("D2").Value = 48.443473941970716
$ws.Range("E2").Value = 53.771243284480605

# Row 3 data values (columns B-E) - STR meanEMG legmaxROM updated
$ws.Range("B3").Value = 71.788163308529889
$ws.Range("C3").Value = 31.805776781676283
$ws.Range("D3").Value = 31.717655404642731
$ws.Range("E3").Value = 56.955511943931079

# Reflect the new active selection range (B1:E3) as in the saved workbook
$ws.Range("B1:E3").Select() | Out-Null
